# Update Visium v2.0 to use old dataset type branch
#
# 1. dataset_type sheet: split "GeoMx" into "GeoMx (NGS)" / "GeoMx (nCounter)",
#    drop "RNAseq (Visium)" and "RNAseq (GeoMx)" (net 35 -> 34 rows).
# 2. Main "Visium" sheet: point the dataset_type validation at the new A1:A34 range.
# 3. .metadata sheet: bump pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

$wsVisium = $wb.Worksheets.Item("Visium")
$wsDatasetType = $wb.Worksheets.Item("dataset_type")
$wsMeta = $wb.Worksheets.Item(".metadata")

# --- 1. Rewrite the dataset_type lookup list -------------------------------

$newDatasetTypes = @(
  @("HiFi-Slide", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000195"),
  @("SNARE-seq2", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000264"),
  @("MIBI", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000172"),
  @("DESI", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000204"),
  @("scATACseq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000247"),
  @("Auto-fluorescence", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000205"),
  @("Confocal", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000206"),
  @("scRNAseq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000248"),
  @("Xenium", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000219"),
  @("snATACseq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000183"),
  @("Molecular Cartography", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000217"),
  @("CosMx", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000218"),
  @("DBiT", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000222"),
  @("SIMS", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000202"),
  @("Cell DIVE", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000159"),
  @("CODEX", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000160"),
  @("GeoMx (NGS)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000300"),
  @("CyCIF", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000200"),
  @("Light Sheet", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000168"),
  @("RNAseq (bulk)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000212"),
  @("MALDI", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000201"),
  @("2D Imaging Mass Cytometry", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296"),
  @("Histology", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000197"),
  @("Enhanced Stimulated Raman Spectroscopy (SRS)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000209"),
  @("ATACseq (bulk)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000210"),
  @("MERFISH", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000221"),
  @("LC-MS", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000194"),
  @("10X Multiome", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000215"),
  @("GeoMx (nCounter)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000301"),
  @("PhenoCycler", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000199"),
  @("Second Harmonic Generation (SHG)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000208"),
  @("Thick section Multiphoton MxIF", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000207"),
  @("snRNAseq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000184"),
  @("Visium", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000187")
)

$rowCount = $newDatasetTypes.Count

# Write the new 34-row table cell-by-cell over the old 35-row one ...
for ($i = 0; $i -lt $rowCount; $i++) {
    $pair = $newDatasetTypes[$i]
    $wsDatasetType.Cells.Item($i + 1, 1).Value = $pair[0]
    $wsDatasetType.Cells.Item($i + 1, 2).Value = $pair[1]
}
# ... then drop the now-stale 35th row entirely.
$wsDatasetType.Rows.Item($rowCount + 1).Delete()

# --- 2. Point the dataset_type validation at the shrunk range --------------

$dv = $wsVisium.Range("D2:D1001").Validation
$dv.Delete()
$dv.Add(3, 1, 1, "'dataset_type'!`$A`$1:`$A`$$rowCount")
$dv.IgnoreBlank = $true
$dv.InCellDropdown = $true
$dv.ErrorTitle = "Validation Error"
$dv.ErrorMessage = ""
$dv.ShowError = $true

# --- 3. Bump the pav:createdOn timestamp on .metadata -----------------------

$wsMeta.Cells.Item(2, 3).Value = "2023-11-15T17:24:29-08:00"
